$d = $word.ActiveDocument

# Update the header date line
$d.Content.Find.Execute("2023-11-16 Thursday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023-11-17 Friday", 2)

$t = $d.Tables(1)

# Each content row in the table holds 5 division problems; blank rows
# separate the groups. Content lives in table rows 1, 5, 9, 13, 17.
$updates = @(
    @{ Row = 1;  Col = 1; Text = "70÷6=11, 4" },
    @{ Row = 1;  Col = 2; Text = "92÷7=13, 1" },
    @{ Row = 1;  Col = 3; Text = "49÷6=8, 1" },
    @{ Row = 1;  Col = 4; Text = "10÷5=2, 0" },
    @{ Row = 1;  Col = 5; Text = "93÷3=31, 0" },

    @{ Row = 5;  Col = 1; Text = "51÷6=8, 3" },
    @{ Row = 5;  Col = 2; Text = "69÷2=34, 1" },
    @{ Row = 5;  Col = 3; Text = "20÷9=2, 2" },
    @{ Row = 5;  Col = 4; Text = "18÷5=3, 3" },
    @{ Row = 5;  Col = 5; Text = "78÷3=26, 0" },

    @{ Row = 9;  Col = 1; Text = "69÷5=13, 4" },
    @{ Row = 9;  Col = 2; Text = "96÷8=12, 0" },
    @{ Row = 9;  Col = 3; Text = "78÷3=26, 0" },
    @{ Row = 9;  Col = 4; Text = "69÷5=13, 4" },
    @{ Row = 9;  Col = 5; Text = "15÷9=1, 6" },

    @{ Row = 13; Col = 1; Text = "14÷8=1, 6" },
    @{ Row = 13; Col = 2; Text = "38÷3=12, 2" },
    @{ Row = 13; Col = 3; Text = "73÷7=10, 3" },
    @{ Row = 13; Col = 4; Text = "50÷5=10, 0" },
    @{ Row = 13; Col = 5; Text = "87÷2=43, 1" },

    @{ Row = 17; Col = 1; Text = "84÷9=9, 3" },
    @{ Row = 17; Col = 2; Text = "12÷3=4, 0" },
    @{ Row = 17; Col = 3; Text = "68÷7=9, 5" },
    @{ Row = 17; Col = 4; Text = "67÷7=9, 4" },
    @{ Row = 17; Col = 5; Text = "19÷6=3, 1" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $r = $cell.Range
    $r.End = $r.End - 1
    $r.Text = $u.Text
}
